# Apply updated price/profit figures pulled by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4830.805
$ws.Range("I138").Value = 5876.857
$ws.Range("J138").Value = 4598.349
$ws.Range("K138").Value = 17630.571
$ws.Range("L138").Value = 13795.047
$ws.Range("M138").Value = -12490.571
$ws.Range("N138").Value = -24075.047

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12059.389
$ws.Range("I32").Value = 10579.97
$ws.Range("J32").Value = 28333
$ws.Range("K32").Value = 10579.97
$ws.Range("L32").Value = 28333
$ws.Range("M32").Value = -10292.97
$ws.Range("N32").Value = -28907

$ws.Range("H74").Value = 2135.3333
$ws.Range("I74").Value = 704
$ws.Range("J74").Value = 3566.6667
$ws.Range("K74").Value = 704
$ws.Range("L74").Value = 3566.6667
$ws.Range("M74").Value = 170
$ws.Range("N74").Value = -5314.6667

$ws.Range("H77").Value = 2135.3333
$ws.Range("I77").Value = 704
$ws.Range("J77").Value = 3566.6667
$ws.Range("K77").Value = 3520
$ws.Range("L77").Value = 17833.3335
$ws.Range("M77").Value = 848
$ws.Range("N77").Value = -26569.3335

$ws.Range("H110").Value = 918.36365
$ws.Range("I110").Value = 480.25
$ws.Range("K110").Value = 480.25
$ws.Range("M110").Value = 1564.75

$ws.Range("H132").Value = 2063.4707
$ws.Range("I132").Value = 1837
$ws.Range("K132").Value = 5511
$ws.Range("M132").Value = -2981

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("N132").ClearContents()

$ws.Range("H134").Value = 7847.75
$ws.Range("I134").Value = 8106.952
$ws.Range("J134").Value = 6033.3335
$ws.Range("K134").Value = 24320.856
$ws.Range("L134").Value = 18100.0005
$ws.Range("M134").Value = -21785.856
$ws.Range("N134").Value = -23170.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2766.5557
$ws.Range("I31").Value = 2362.5
$ws.Range("J31").Value = 5999
$ws.Range("K31").Value = 2362.5
$ws.Range("L31").Value = 5999
$ws.Range("M31").Value = -2067.5
$ws.Range("N31").Value = -6589

$ws.Range("H34").Value = 2766.5557
$ws.Range("I34").Value = 2362.5
$ws.Range("J34").Value = 5999
$ws.Range("K34").Value = 2362.5
$ws.Range("L34").Value = 5999
$ws.Range("M34").Value = -2160.5
$ws.Range("N34").Value = -6403

$ws.Range("H58").Value = 2290367.5
$ws.Range("I58").Value = 3107338.8
$ws.Range("J58").Value = 2848.2
$ws.Range("K58").Value = 3107338.8
$ws.Range("L58").Value = 2848.2
$ws.Range("M58").Value = -3107135.8
$ws.Range("N58").Value = -3254.2

$ws.Range("H132").Value = 2367.5833
$ws.Range("I132").Value = 1641.1
$ws.Range("K132").Value = 4923.299999999999
$ws.Range("M132").Value = -2393.299999999999

$ws.Range("H134").Value = 1165.5238
$ws.Range("I134").Value = 1041.8438
$ws.Range("K134").Value = 3125.5314
$ws.Range("M134").Value = -590.5314000000003

$ws.Range("H136").Value = 2290367.5
$ws.Range("I136").Value = 3107338.8
$ws.Range("J136").Value = 2848.2
$ws.Range("K136").Value = 9322016.399999999
$ws.Range("L136").Value = 8544.599999999999
$ws.Range("M136").Value = -9319466.399999999
$ws.Range("N136").Value = -13644.6

$ws.Range("H141").Value = 194926
$ws.Range("J141").Value = 194926
$ws.Range("L141").Value = 194926
$ws.Range("N141").Value = -205286

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 527448.0600000001
$ws.Range("I4").Value = 500416.8
$ws.Range("K4").Value = 1501250.4
$ws.Range("M4").Value = -1501138.4

$ws.Range("H34").Value = 1356.909
$ws.Range("J34").Value = 2432.6667
$ws.Range("L34").Value = 7298.000100000001
$ws.Range("N34").Value = -7466.000100000001

$ws.Range("H113").Value = 80285.57000000001
$ws.Range("J113").Value = 1999.5
$ws.Range("L113").Value = 5998.5
$ws.Range("N113").Value = -10338.5

$ws.Range("H122").Value = 1477.6666
$ws.Range("J122").Value = 1949.5
$ws.Range("L122").Value = 17545.5
$ws.Range("N122").Value = -22445.5

$ws.Range("H131").Value = 11256.921
$ws.Range("J131").Value = 11696.795
$ws.Range("L131").Value = 35090.385
$ws.Range("N131").Value = -45170.385

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 19333.334
$ws.Range("J46").Value = 19333.334
$ws.Range("L46").Value = 19333.334
$ws.Range("N46").Value = -19645.334

$ws.Range("H102").Value = 2458.5334
$ws.Range("I102").Value = 2243.4546
$ws.Range("K102").Value = 2243.4546
$ws.Range("M102").Value = -621.4546

$ws.Range("H107").Value = 470.35715
$ws.Range("I107").Value = 298.375
$ws.Range("K107").Value = 298.375
$ws.Range("M107").Value = 1621.625

$ws.Range("H132").Value = 1610322.5
$ws.Range("I132").Value = 2573797.8
$ws.Range("K132").Value = 7721393.399999999
$ws.Range("M132").Value = -7718863.399999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6511.8
$ws.Range("J7").Value = 2500
$ws.Range("L7").Value = 2500
$ws.Range("N7").Value = -2724

$ws.Range("H40").Value = 16506.312
$ws.Range("I40").Value = 23033.666
$ws.Range("J40").Value = 12589.9
$ws.Range("K40").Value = 23033.666
$ws.Range("L40").Value = 12589.9
$ws.Range("M40").Value = -22897.666
$ws.Range("N40").Value = -12861.9

$ws.Range("H126").Value = 6511.8
$ws.Range("J126").Value = 2500
$ws.Range("L126").Value = 7500
$ws.Range("N126").Value = -12440

$ws.Range("H136").Value = 4527.5557
$ws.Range("I136").Value = 3549.6
$ws.Range("J136").Value = 5750
$ws.Range("K136").Value = 10648.8
$ws.Range("L136").Value = 17250
$ws.Range("M136").Value = -8098.799999999999
$ws.Range("N136").Value = -22350

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1498.4445
$ws.Range("I81").Value = 1398.25
$ws.Range("K81").Value = 2796.5
$ws.Range("M81").Value = -1735.5

$ws.Range("H84").Value = 1498.4445
$ws.Range("I84").Value = 1398.25
$ws.Range("K84").Value = 13982.5
$ws.Range("M84").Value = -8678.5

$ws.Range("H107").Value = 895.625
$ws.Range("I107").Value = 743.6667
$ws.Range("K107").Value = 2231.0001
$ws.Range("M107").Value = -311.0001000000002

$ws.Range("H126").Value = 3969.7856
$ws.Range("I126").Value = 5151.4287
$ws.Range("K126").Value = 15454.2861
$ws.Range("M126").Value = -12984.2861

$ws.Range("H132").Value = 1480.9375
$ws.Range("I132").Value = 1053.1538
$ws.Range("J132").Value = 3334.6667
$ws.Range("K132").Value = 3159.4614
$ws.Range("L132").Value = 10004.0001
$ws.Range("M132").Value = -629.4614000000001
$ws.Range("N132").Value = -15064.0001
